# Data Resume_Camryn Allen.docx edit:
#   - Remove the reviewer comment anchored on "ON" (the back half of the
#     "EDUCATION" heading), which in turn collapses the two split runs
#     ("EDUCATI" + "ON") that existed only because the comment range split
#     them back into a single "EDUCATION" run.

$d = $word.ActiveDocument

# Delete the (only) comment in the document. This removes the
# commentRangeStart/commentRangeEnd/commentReference markers around "ON"
# along with the comment content itself.
if ($d.Comments.Count -gt 0) {
    for ($i = $d.Comments.Count; $i -ge 1; $i--) {
        $d.Comments($i).Delete()
    }
}

# Re-run the heading text through Find/Replace so the now-adjacent
# "EDUCATI" + "ON" runs are coalesced back into a single "EDUCATION" run,
# matching how Word normalizes the run after the comment boundaries that
# used to separate them are gone.
$d.Content.Find.Execute("EDUCATION", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "EDUCATION", 2) | Out-Null
